# Update threat_modeling.xlsx:
#  - sheet "threat_list": fix wording of the password-policy mitigation text
#    (used in H5, H9, H24), and normalise the formatting of rows 28-29
#    (TR-62 / TR-63) to match the rest of the table (row 27's style).
#  - sheet "threat_modeling_all": add the two missing threats (TR-62, TR-63)
#    as new rows 65-66, mirroring columns A-F of threat_list rows 28-29.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("threat_modeling_all")
$ws2 = $wb.Worksheets.Item("threat_list")

# ---------------------------------------------------------------------
# 1) threat_list: correct the cryptographically-strong-password wording
# ---------------------------------------------------------------------
$newPasswordText = "Apply setting policy of cryptographically strong password`n- Enforce passwords longer than 7 characters.`n- Forces the use of mixed the letters of the alphabet and numbers.`nStrong authentication method`n- Condider 2-Factor-Authentication method"

$ws2.Cells.Item(5, 8).Value = $newPasswordText
$ws2.Cells.Item(9, 8).Value = $newPasswordText
$ws2.Cells.Item(24, 8).Value = $newPasswordText

# ---------------------------------------------------------------------
# 2) threat_list: re-format rows 28 and 29 (TR-62 / TR-63) so they use
#    the same cell styles as the rest of the table (like row 27) instead
#    of the old one-off "green highlight" styles.
# ---------------------------------------------------------------------
function Copy-RowFormat($srcRow, $dstRow) {
    $savedValues = @()
    for ($c = 1; $c -le 8; $c++) {
        $savedValues += , $ws2.Cells.Item($dstRow, $c).Value2
    }

    $ws2.Range("A$srcRow`:H$srcRow").Copy() | Out-Null
    $ws2.Range("A$dstRow`:H$dstRow").PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = $false

    for ($c = 1; $c -le 8; $c++) {
        $ws2.Cells.Item($dstRow, $c).Value = $savedValues[$c - 1]
    }
}

Copy-RowFormat 27 28
Copy-RowFormat 27 29

# ---------------------------------------------------------------------
# 3) threat_modeling_all: append TR-62 and TR-63 as rows 65 and 66,
#    copying the row-64 style and filling in columns A-F only.
# ---------------------------------------------------------------------
$ws1.Range("A64:H64").Copy() | Out-Null
$ws1.Range("A65:H65").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws1.Range("A64:H64").Copy() | Out-Null
$ws1.Range("A66:H66").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 65 <- threat_list row 28 (TR-62)
$ws1.Cells.Item(65, 1).Value = $ws2.Cells.Item(28, 1).Value2
$ws1.Cells.Item(65, 2).Value = $ws2.Cells.Item(28, 2).Value2
$ws1.Cells.Item(65, 3).Value = $ws2.Cells.Item(28, 3).Value2
$ws1.Cells.Item(65, 4).Value = $ws2.Cells.Item(28, 4).Value2
$ws1.Cells.Item(65, 5).Value = $ws2.Cells.Item(28, 5).Value2
$ws1.Cells.Item(65, 6).Value = $ws2.Cells.Item(28, 6).Value2
$ws1.Cells.Item(65, 7).Value = ""
$ws1.Cells.Item(65, 8).Value = ""
$ws1.Rows.Item(65).RowHeight = 66

# Row 66 <- threat_list row 29 (TR-63)
$ws1.Cells.Item(66, 1).Value = $ws2.Cells.Item(29, 1).Value2
$ws1.Cells.Item(66, 2).Value = $ws2.Cells.Item(29, 2).Value2
$ws1.Cells.Item(66, 3).Value = $ws2.Cells.Item(29, 3).Value2
$ws1.Cells.Item(66, 4).Value = $ws2.Cells.Item(29, 4).Value2
$ws1.Cells.Item(66, 5).Value = $ws2.Cells.Item(29, 5).Value2
$ws1.Cells.Item(66, 6).Value = $ws2.Cells.Item(29, 6).Value2
$ws1.Cells.Item(66, 7).Value = ""
$ws1.Cells.Item(66, 8).Value = ""
$ws1.Rows.Item(66).RowHeight = 82.5

# ---------------------------------------------------------------------
# 4) Views / selections
# ---------------------------------------------------------------------
$ws1.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 56
$win.ScrollColumn = 1
$ws1.Range("A64:H66").Select() | Out-Null

$ws2.Activate()
$ws2.Range("H26").Select() | Out-Null

Write-Host "edit complete"
